$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Centraal Station_A"
$ws.Range("A18").Select()
